$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column: header in H1 styled like the existing header row (copy
# format from G1, the "sum" header immediately to its left), data rows
# H2:H4 all zero.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
